# Auto-generated Excel COM-interop edit script
# Applies the Sagittarius_Profits market-data refresh (scheduled runner update)
# to columns H-N (currentAveragePrice*, LevePrice*, LeveProfit*) across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2126.0527
$ws.Range("I40").Value = 2061.8333
$ws.Range("J40").Value = 2236.1428
$ws.Range("K40").Value = 2061.8333
$ws.Range("L40").Value = 2236.1428
$ws.Range("M40").Value = -1886.8333
$ws.Range("N40").Value = -2586.1428

$ws.Range("H64").Value = 5239.8
$ws.Range("I64").Value = 4850
$ws.Range("K64").Value = 4850
$ws.Range("M64").Value = -4602

$ws.Range("H67").Value = 5239.8
$ws.Range("I67").Value = 4850
$ws.Range("K67").Value = 4850
$ws.Range("M67").Value = -3992

$ws.Range("H88").Value = 5212.615
$ws.Range("I88").Value = 690.1667
$ws.Range("J88").Value = 9089
$ws.Range("K88").Value = 690.1667
$ws.Range("L88").Value = 9089
$ws.Range("M88").Value = -284.1667
$ws.Range("N88").Value = -9901

$ws.Range("H91").Value = 5212.615
$ws.Range("I91").Value = 690.1667
$ws.Range("J91").Value = 9089
$ws.Range("K91").Value = 690.1667
$ws.Range("L91").Value = 9089
$ws.Range("M91").Value = 713.8333
$ws.Range("N91").Value = -11897

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6189.5713
$ws.Range("I45").Value = 5259.8
$ws.Range("K45").Value = 5259.8
$ws.Range("M45").Value = -4882.8

$ws.Range("H88").Value = 766.3333
$ws.Range("I88").Value = 799
$ws.Range("K88").Value = 799
$ws.Range("M88").Value = -393

$ws.Range("H91").Value = 766.3333
$ws.Range("I91").Value = 799
$ws.Range("K91").Value = 799
$ws.Range("M91").Value = 605

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 2900
$ws.Range("I37").Value = 2900
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 2900
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -2763
$ws.Range("N37").ClearContents()

$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()

$ws.Range("H86").Value = 7775
$ws.Range("I86").Value = 7033.3335
$ws.Range("J86").Value = 10000
$ws.Range("K86").Value = 7033.3335
$ws.Range("L86").Value = 10000
$ws.Range("M86").Value = -5910.3335
$ws.Range("N86").Value = -12246

$ws.Range("H89").Value = 7775
$ws.Range("I89").Value = 7033.3335
$ws.Range("J89").Value = 10000
$ws.Range("K89").Value = 35166.6675
$ws.Range("L89").Value = 50000
$ws.Range("M89").Value = -29550.6675
$ws.Range("N89").Value = -61232

$ws.Range("H94").Value = 1701.7273
$ws.Range("I94").Value = 1722.4
$ws.Range("J94").Value = 1495
$ws.Range("K94").Value = 1722.4
$ws.Range("L94").Value = 1495
$ws.Range("M94").Value = -1271.4
$ws.Range("N94").Value = -2397

$ws.Range("H99").Value = 2013.9375
$ws.Range("I99").Value = 2022.8667
$ws.Range("K99").Value = 2022.8667
$ws.Range("M99").Value = -524.8667

$ws.Range("H105").Value = 3998
$ws.Range("I105").Value = 3998
$ws.Range("K105").Value = 3998
$ws.Range("M105").Value = -2251

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2222
$ws.Range("I31").Value = 1844.625
$ws.Range("K31").Value = 1844.625
$ws.Range("M31").Value = -1549.625

$ws.Range("H34").Value = 2222
$ws.Range("I34").Value = 1844.625
$ws.Range("K34").Value = 1844.625
$ws.Range("M34").Value = -1642.625

$ws.Range("H58").Value = 2127.6316
$ws.Range("I58").Value = 2084.1177
$ws.Range("K58").Value = 2084.1177
$ws.Range("M58").Value = -1881.1177

$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws.Range("H136").Value = 2127.6316
$ws.Range("I136").Value = 2084.1177
$ws.Range("K136").Value = 6252.353099999999
$ws.Range("M136").Value = -3702.353099999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 388.75
$ws.Range("I44").Value = 303
$ws.Range("J44").Value = 474.5
$ws.Range("K44").Value = 909
$ws.Range("L44").Value = 1423.5
$ws.Range("M44").Value = -511
$ws.Range("N44").Value = -2219.5

$ws.Range("H113").Value = 996.3333
$ws.Range("I113").Value = 553.2857
$ws.Range("J113").Value = 1616.6
$ws.Range("K113").Value = 1659.8571
$ws.Range("L113").Value = 4849.799999999999
$ws.Range("M113").Value = 510.1428999999998
$ws.Range("N113").Value = -9189.799999999999

$ws.Range("H114").Value = 1968.8823
$ws.Range("J114").Value = 1854.1
$ws.Range("L114").Value = 5562.299999999999
$ws.Range("N114").Value = -12070.3

$ws.Range("H129").Value = 2468.3845
$ws.Range("J129").Value = 3018.7778
$ws.Range("L129").Value = 9056.3334
$ws.Range("N129").Value = -19056.3334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 119.85714
$ws.Range("I2").Value = 147.36363
$ws.Range("K2").Value = 147.36363
$ws.Range("M2").Value = -34.36363

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()

$ws.Range("H47").Value = 35031
$ws.Range("J47").Value = 35031
$ws.Range("L47").Value = 35031
$ws.Range("N47").Value = -36167

$ws.Range("H70").Value = 8113.467
$ws.Range("I70").Value = 7988.375
$ws.Range("J70").Value = 8256.429
$ws.Range("K70").Value = 7988.375
$ws.Range("L70").Value = 8256.429
$ws.Range("M70").Value = -7718.375
$ws.Range("N70").Value = -8796.429

$ws.Range("H73").Value = 8113.467
$ws.Range("I73").Value = 7988.375
$ws.Range("J73").Value = 8256.429
$ws.Range("K73").Value = 7988.375
$ws.Range("L73").Value = 8256.429
$ws.Range("M73").Value = -7052.375
$ws.Range("N73").Value = -10128.429

$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3751.9
$ws.Range("I40").Value = 2305
$ws.Range("K40").Value = 2305
$ws.Range("M40").Value = -2169

$ws.Range("H46").Value = 47688.363
$ws.Range("J46").Value = 2500
$ws.Range("L46").Value = 2500
$ws.Range("N46").Value = -2876

$ws.Range("H58").Value = 693
$ws.Range("I58").Value = 693
$ws.Range("K58").Value = 693
$ws.Range("M58").Value = -433

$ws.Range("H61").Value = 4072.6667
$ws.Range("I61").Value = 5402.3335
$ws.Range("K61").Value = 5402.3335
$ws.Range("M61").Value = -5200.3335

$ws.Range("H113").Value = 4072.6667
$ws.Range("I113").Value = 5402.3335
$ws.Range("K113").Value = 5402.3335
$ws.Range("M113").Value = -3232.3335

$ws.Range("H122").Value = 8117.909
$ws.Range("J122").Value = 6729.2
$ws.Range("L122").Value = 20187.6
$ws.Range("N122").Value = -25087.6

$ws.Range("H132").Value = 2999.1
$ws.Range("I132").Value = 2398.25
$ws.Range("K132").Value = 7194.75
$ws.Range("M132").Value = -4664.75

$ws.Range("H136").Value = 2862.1428
$ws.Range("I136").Value = 2364.4285
$ws.Range("J136").Value = 3857.5715
$ws.Range("K136").Value = 7093.2855
$ws.Range("L136").Value = 11572.7145
$ws.Range("M136").Value = -4543.2855
$ws.Range("N136").Value = -16672.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 501.25
$ws.Range("I14").Value = 152.5
$ws.Range("J14").Value = 850
$ws.Range("K14").Value = 152.5
$ws.Range("L14").Value = 850
$ws.Range("M14").Value = 15.5
$ws.Range("N14").Value = -1186

$ws.Range("H62").Value = 20851.111
$ws.Range("I62").Value = 32932.668
$ws.Range("J62").Value = 14810.333
$ws.Range("K62").Value = 32932.668
$ws.Range("L62").Value = 14810.333
$ws.Range("M62").Value = -32308.668
$ws.Range("N62").Value = -16058.333

$ws.Range("H65").Value = 20851.111
$ws.Range("I65").Value = 32932.668
$ws.Range("J65").Value = 14810.333
$ws.Range("K65").Value = 164663.34
$ws.Range("L65").Value = 74051.66500000001
$ws.Range("M65").Value = -161543.34
$ws.Range("N65").Value = -80291.66500000001

$ws.Range("H81").Value = 1667232
$ws.Range("I81").Value = 660.25
$ws.Range("K81").Value = 1320.5
$ws.Range("M81").Value = -259.5

$ws.Range("H84").Value = 1667232
$ws.Range("I84").Value = 660.25
$ws.Range("K84").Value = 6602.5
$ws.Range("M84").Value = -1298.5

$ws.Range("H132").Value = 4324.5
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws.Range("H136").Value = 2243.4707
$ws.Range("I136").Value = 1876
$ws.Range("J136").Value = 4999.5
$ws.Range("K136").Value = 5628
$ws.Range("L136").Value = 14998.5
$ws.Range("M136").Value = -3078
$ws.Range("N136").Value = -20098.5
